$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B15").Value = "29/03/2022 11:51:40"
$ws.Range("C15").Value = "29/03/2022 11:51:43"
$ws.Range("D15").Value = "00:00:02"
$ws.Range("E15").Value = "4"
$ws.Range("A16").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B16").Value = "29/03/2022 12:15:05"
$ws.Range("C16").Value = "29/03/2022 12:15:08"
$ws.Range("D16").Value = "00:00:02"
$ws.Range("E16").Value = "4"
$ws.Range("A17").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B17").Value = "29/03/2022 12:36:09"
$ws.Range("C17").Value = "29/03/2022 12:36:12"
$ws.Range("D17").Value = "00:00:03"
$ws.Range("E17").Value = "7"
$ws.Range("A18").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B18").Value = "29/03/2022 12:56:55"
$ws.Range("C18").Value = "29/03/2022 12:56:55"
$ws.Range("D18").Value = "00:00:00"
$ws.Range("E18").Value = "1"
$ws.Range("A19").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B19").Value = "29/03/2022 13:21:08"
$ws.Range("C19").Value = "29/03/2022 13:24:49"
$ws.Range("D19").Value = "00:03:41"
$ws.Range("E19").Value = "2"
$ws.Range("A20").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B20").Value = "29/03/2022 13:30:40"
$ws.Range("C20").Value = "29/03/2022 13:32:13"
$ws.Range("D20").Value = "00:01:32"
$ws.Range("E20").Value = "3"
$ws.Range("A21").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B21").Value = "29/03/2022 13:38:19"
$ws.Range("C21").Value = "29/03/2022 13:42:06"
$ws.Range("D21").Value = "00:03:46"
$ws.Range("E21").Value = "2"
$ws.Range("A22").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B22").Value = "29/03/2022 15:44:24"
$ws.Range("C22").Value = "29/03/2022 15:55:50"
$ws.Range("D22").Value = "00:11:26"
$ws.Range("E22").Value = "6"
$ws.Range("A23").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B23").Value = "29/03/2022 16:25:42"
$ws.Range("C23").Value = "29/03/2022 16:27:40"
$ws.Range("D23").Value = "00:01:58"
$ws.Range("E23").Value = "1"
$ws.Range("A24").Value = "prj_HC_ContasMunicipais_SANASA"
$ws.Range("B24").Value = "29/03/2022 17:12:56"
$ws.Range("C24").Value = "29/03/2022 17:20:10"
$ws.Range("D24").Value = "00:07:13"
$ws.Range("E24").Value = "5"
